# Apply the "FilesTab" query update (the text in Sheet1!B5):
#  - remove the UNION with df_sequencing_file, keep only df_pathology_file
#  - reformat the File Size CASE expression to use
#      RTRIM(RTRIM(printf('%.2f', ...), '0'), '.')  instead of ROUND(..., 2)
#  - split the WHERE clause across multiple AND lines
#  - add "Order by fd.file_name asc" before the LIMIT

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newQuery = "with file_data as (select file_name, data_category,file_type, file_size,file_access,  file_description,`"sample.id`" from df_pathology_file)`nSELECT DISTINCT`n    fd.file_name AS `"File Name`",`n    fd.data_category AS `"Data Category`",`n    COALESCE(fd.file_description, '') AS `"File Description`",`n    fd.file_type AS `"File Type`",`n    CASE     `n        WHEN fd.file_size >= 1024 * 1024 * 1024 THEN `n            RTRIM(RTRIM(printf('%.2f', fd.file_size / (1024.0 * 1024.0 * 1024.0)), '0'), '.') || ' GB'`n        WHEN fd.file_size >= 1024 * 1024 THEN `n            RTRIM(RTRIM(printf('%.2f', fd.file_size / (1024.0 * 1024.0)), '0'), '.') || ' MB'`n        WHEN fd.file_size >= 1024 THEN `n            RTRIM(RTRIM(printf('%.2f', fd.file_size / 1024.0), '0'), '.') || ' KB'`n        ELSE `n            RTRIM(RTRIM(printf('%.2f', fd.file_size), '0'), '.') || ' Bytes'`n    END AS `"File Size`",`n    fd.file_access AS `"File Access`",`n    std.dbgap_accession AS `"Study ID`",`n    prt.participant_id AS `"Participant ID`",`n    smp.sample_id AS `"Sample ID`"    `nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_sample smp ON prt.id = smp.`"participant.id`"`nJOIN `n    file_data fd ON smp.id = fd.`"sample.id`"`nWHERE `n    std.dbgap_accession = 'phs000720' `n    AND fd.file_type = 'dicom' `n    AND smp.tumor_classification = 'Metastatic'`nOrder by fd.file_name asc LIMIT 100;"

$ws.Range("B5").Value = $newQuery

# The author's Excel session also left the active selection on B5
# (was A5 before), so move the selection there to match.
$ws.Range("B5").Select()
